$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (E1) into the new
# header cell F1, so the new "time_taken" header gets the same bold,
# centered, bordered style used by the other header cells.
$ws.Range("E1").Copy($ws.Range("F1"))

# Set the header text for the new column.
$ws.Range("F1").Value = "time_taken"

# Add the data value for the new column (plain data cell, no special style,
# matching the unstyled data cells in row 2).
$ws.Range("F2").Value = "2021-10-05 13:41:40.871703"
